$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 21903.771
$ws.Range("I135").Value = 27902.73
$ws.Range("J135").Value = 1725.4546
$ws.Range("K135").Value = 251124.57
$ws.Range("L135").Value = 15529.0914
$ws.Range("M135").Value = -248589.57
$ws.Range("N135").Value = -20599.0914

$ws.Range("H137").Value = 4688949
$ws.Range("I137").Value = 2084613.9
$ws.Range("J137").Value = 12501954
$ws.Range("K137").Value = 6253841.699999999
$ws.Range("L137").Value = 37505862
$ws.Range("M137").Value = -6251291.699999999
$ws.Range("N137").Value = -37510962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 15183

$ws.Range("H55").Value = 16049.667
$ws.Range("J55").Value = 16049.667
$ws.Range("L55").Value = 16049.667
$ws.Range("N55").Value = -16679.667

$ws.Range("H61").Value = 1479.5
$ws.Range("I61").Value = 1523.7084
$ws.Range("J61").Value = 1302.6666
$ws.Range("K61").Value = 1523.7084
$ws.Range("L61").Value = 1302.6666
$ws.Range("M61").Value = -1311.7084
$ws.Range("N61").Value = -1726.6666

$ws.Range("H80").Value = 24661.5
$ws.Range("J80").Value = 24661.5
$ws.Range("L80").Value = 24661.5
$ws.Range("N80").Value = -26657.5

$ws.Range("H83").Value = 24661.5
$ws.Range("J83").Value = 24661.5
$ws.Range("L83").Value = 73984.5
$ws.Range("N83").Value = -83968.5

$ws.Range("H121").Value = 31542.5
$ws.Range("J121").Value = 31542.5
$ws.Range("L121").Value = 31542.5
$ws.Range("N121").Value = -35036.5

$ws.Range("H136").Value = 1479.5
$ws.Range("I136").Value = 1523.7084
$ws.Range("J136").Value = 1302.6666
$ws.Range("K136").Value = 4571.1252
$ws.Range("L136").Value = 3907.9998
$ws.Range("M136").Value = -2021.1252
$ws.Range("N136").Value = -9007.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34987
$ws.Range("J35").Value = 34987
$ws.Range("L35").Value = 34987
$ws.Range("N35").Value = -35607

$ws.Range("H82").Value = 20413.223
$ws.Range("I82").Value = 10257
$ws.Range("J82").Value = 21682.75
$ws.Range("K82").Value = 10257
$ws.Range("L82").Value = 21682.75
$ws.Range("M82").Value = -9874
$ws.Range("N82").Value = -22448.75

$ws.Range("H85").Value = 20413.223
$ws.Range("I85").Value = 10257
$ws.Range("J85").Value = 21682.75
$ws.Range("K85").Value = 10257
$ws.Range("L85").Value = 21682.75
$ws.Range("M85").Value = -8931
$ws.Range("N85").Value = -24334.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2181.9092
$ws.Range("I31").Value = 1774.2667
$ws.Range("J31").Value = 3055.4285
$ws.Range("K31").Value = 1774.2667
$ws.Range("L31").Value = 3055.4285
$ws.Range("M31").Value = -1479.2667
$ws.Range("N31").Value = -3645.4285

$ws.Range("H34").Value = 2181.9092
$ws.Range("I34").Value = 1774.2667
$ws.Range("J34").Value = 3055.4285
$ws.Range("K34").Value = 1774.2667
$ws.Range("L34").Value = 3055.4285
$ws.Range("M34").Value = -1572.2667
$ws.Range("N34").Value = -3459.4285

$ws.Range("H41").Value = 16232.5
$ws.Range("J41").Value = 19976.666
$ws.Range("L41").Value = 19976.666
$ws.Range("N41").Value = -20832.666

$ws.Range("H50").Value = 8691.857
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 9307.166999999999
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 9307.166999999999
$ws.Range("M50").Value = -4375
$ws.Range("N50").Value = -10557.167

$ws.Range("H51").Value = 8749
$ws.Range("I51").Value = 4500
$ws.Range("J51").Value = 9598.799999999999
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 9598.799999999999
$ws.Range("M51").Value = -3764
$ws.Range("N51").Value = -11070.8

$ws.Range("H60").Value = 142124
$ws.Range("I60").Value = 4000
$ws.Range("J60").Value = 161856
$ws.Range("K60").Value = 4000
$ws.Range("L60").Value = 161856
$ws.Range("M60").Value = -3489
$ws.Range("N60").Value = -162878

$ws.Range("H61").Value = 8749
$ws.Range("I61").Value = 4500
$ws.Range("J61").Value = 9598.799999999999
$ws.Range("K61").Value = 4500
$ws.Range("L61").Value = 9598.799999999999
$ws.Range("M61").Value = -4152
$ws.Range("N61").Value = -10294.8

$ws.Range("H68").Value = 18295
$ws.Range("J68").Value = 18295
$ws.Range("L68").Value = 18295
$ws.Range("N68").Value = -19793

$ws.Range("H71").Value = 18295
$ws.Range("J71").Value = 18295
$ws.Range("L71").Value = 54885
$ws.Range("N71").Value = -62373

$ws.Range("H99").Value = 1859.2632
$ws.Range("I99").Value = 1692.6666
$ws.Range("J99").Value = 2144.8572
$ws.Range("K99").Value = 1692.6666
$ws.Range("L99").Value = 2144.8572
$ws.Range("M99").Value = -194.6666
$ws.Range("N99").Value = -5140.8572

$ws.Range("H109").Value = 10975
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 10975
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 10975
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -13055

$ws.Range("H126").Value = 1859.2632
$ws.Range("I126").Value = 1692.6666
$ws.Range("J126").Value = 2144.8572
$ws.Range("K126").Value = 5077.9998
$ws.Range("L126").Value = 6434.571599999999
$ws.Range("M126").Value = -2607.9998
$ws.Range("N126").Value = -11374.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 675.119
$ws.Range("J113").Value = 715
$ws.Range("L113").Value = 2145
$ws.Range("N113").Value = -6485

$ws.Range("H131").Value = 972.6429000000001
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 972.14545
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 2916.43635
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -12996.43635

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25061
$ws.Range("J57").Value = 25061
$ws.Range("L57").Value = 25061
$ws.Range("N57").Value = -26701

$ws.Range("H123").Value = 33217.75
$ws.Range("J123").Value = 33217.75
$ws.Range("L123").Value = 33217.75
$ws.Range("N123").Value = -38117.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 19442.5
$ws.Range("J109").Value = 19442.5
$ws.Range("L109").Value = 19442.5
$ws.Range("N109").Value = -22216.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 23194.25
$ws.Range("J109").Value = 23194.25
$ws.Range("L109").Value = 23194.25
$ws.Range("N109").Value = -25968.25
